# Auto-generated Excel COM-interop edit script
# Updates the cryptos list per the commit "Updated cryptos list on Thu May 18 22:25:51 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.899.69"
$ws.Range("E2").Value = "  -1.71%  "
$ws.Range("D3").Value = "1.811.88"
$ws.Range("E3").Value = "  -0.69%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("E7").Value = "  +3.76%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3755"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07469"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8841"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.47"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.58%  "
$ws.Range("D12").Value = "1.821.24"
$ws.Range("E12").Value = "  -0.20%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.362"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.545"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07051"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.47"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.61%  "
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008773"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9999"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.11%  "
$ws.Range("E20").Value = "  -2.99%  "
$ws.Range("D21").Value = "26.893.61"
$ws.Range("E21").Value = "  -1.77%  "
$ws.Range("E22").Value = "  +1.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.82"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.84%  "
$ws.Range("D24").Value = "1.979.50"
$ws.Range("E24").Value = "  -3.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.922"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.62"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.51"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.159"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -9.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.298"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.39"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08895"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.95%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7713"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.68%  "
$ws.Range("E33").Value = "  -2.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.488"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.898"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9997"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("E37").Value = "  +0.50%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.461"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01961"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05243"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.75%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5336"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.39%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.226"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.24%  "
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.909"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.85%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1663"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.73%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.616"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.73%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5079"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.41"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "104.52"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.18%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.675"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.51%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.9995"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06329"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.80%  "
